$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list1")

# Set the value of D6 (previously empty) to "gfd"
$ws.Range("D6").Value = "gfd"

# Update the active selection to D6, matching the recorded selection state
$ws.Range("D6").Select()
